$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.256000000000002
$ws.Range("B3").Value = 7.007000000000001
$ws.Range("D4").Value = -7.479000000000001
$ws.Range("B5").Value = 6.314
$ws.Range("D6").Value = -7.834999999999999
$ws.Range("C7").Value = -13.061
$ws.Range("A9").Value = -21.387
$ws.Range("C9").Value = -12.626
$ws.Range("D10").Value = -7.644000000000001
$ws.Range("B11").Value = 7.122
$ws.Range("D11").Value = -8.454000000000001
$ws.Range("B12").Value = 6.404999999999999
$ws.Range("E12").Value = 13.323
$ws.Range("A13").Value = -21.918
$ws.Range("A16").Value = -20.84
$ws.Range("E17").Value = 13.387
$ws.Range("A18").Value = -21.751
$ws.Range("E19").Value = 13.169
$ws.Range("A20").Value = -21.664
$ws.Range("B21").Value = 6.601999999999999
$ws.Range("C21").Value = -12.282
$ws.Range("D21").Value = -7.747
$ws.Range("E24").Value = 13.454
$ws.Range("D25").Value = -8.039999999999999
